$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "D2" '92.631.92'
Set-TextValue "E2" '  +5.49%  '

Set-TextValue "D3" '3.278.99'
Set-TextValue "E3" '  +0.64%  '

Set-TextValue "E4" '  -0.09%  '

Set-TextValue "D5" '218.52'
Set-TextValue "E5" '  +3.29%  '

Set-TextValue "D6" '631.85'
Set-TextValue "E6" '  +0.86%  '

Set-TextValue "D7" '0.403'
Set-TextValue "E7" '  +5.01%  '

Set-TextValue "E8" '  +3.55%  '

Set-TextValue "D9" '0.998'
Set-TextValue "E9" '  -0.05%  '

Set-TextValue "D10" '3.266.00'
Set-TextValue "E10" '  +0.36%  '

Set-TextValue "D11" '0.592'
Set-TextValue "E11" '  +3.10%  '

Set-TextValue "D12" '0.0000269'
Set-TextValue "E12" '  +2.97%  '

Set-TextValue "E13" '  -2.18%  '

Set-TextValue "D14" '34.42'
Set-TextValue "E14" '  +1.32%  '

Set-TextValue "B15" 'WrappedBTC'
Set-TextValue "C15" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D15" '92.315.81'
Set-TextValue "E15" '  +5.42%  '

Set-TextValue "B16" 'WrappedliquidstakedEther2.0'
Set-TextValue "C16" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D16" '3.878.04'
Set-TextValue "E16" '  +0.64%  '

Set-TextValue "E17" '  +0.28%  '

Set-TextValue "D18" '3.252.07'
Set-TextValue "E18" '  -0.11%  '

Set-TextValue "D19" '3.31'
Set-TextValue "E19" '  +6.39%  '

Set-TextValue "B20" 'PEPE'
Set-TextValue "C20" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D20" '0.0000213'
Set-TextValue "E20" '  +61.32%  '

Set-TextValue "B21" 'Chainlink'
Set-TextValue "C21" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D21" '14.01'
Set-TextValue "E21" '  -0.02%  '

Set-TextValue "B22" 'BitcoinCash'
Set-TextValue "C22" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D22" '445.76'
Set-TextValue "E22" '  +3.03%  '

Set-TextValue "D23" '8.90'
Set-TextValue "E23" '  +0.35%  '

Set-TextValue "D24" '5.28'
Set-TextValue "E24" '  -1.60%  '

Set-TextValue "E25" '  +4.54%  '

Set-TextValue "D26" '12.22'
Set-TextValue "E26" '  -1.06%  '

Set-TextValue "D27" '3.450.19'
Set-TextValue "E27" '  +1.99%  '

Set-TextValue "D28" '77.49'
Set-TextValue "E28" '  +1.49%  '

Set-TextValue "E29" '  -0.03%  '

Set-TextValue "D30" '0.175'
Set-TextValue "E30" '  -2.42%  '

Set-TextValue "D31" '0.998'
Set-TextValue "E31" '  -0.03%  '

Set-TextValue "D32" '8.78'
Set-TextValue "E32" '  -0.14%  '

Set-TextValue "D33" '558.83'
Set-TextValue "E33" '  -0.39%  '

Set-TextValue "B34" 'dogwifhat'
Set-TextValue "C34" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D34" '3.85'
Set-TextValue "E34" '  +30.34%  '

Set-TextValue "B35" 'RenderToken'
Set-TextValue "C35" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue "D35" '7.21'
Set-TextValue "E35" '  +1.06%  '

Set-TextValue "E36" '  -1.10%  '

Set-TextValue "E37" '  -8.03%  '

Set-TextValue "E38" '  +1.05%  '

Set-TextValue "B39" 'WhiteBITCoin'
Set-TextValue "C39" 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D39" '22.44'
Set-TextValue "E39" '  +3.28%  '

Set-TextValue "B40" 'Kaspa'
Set-TextValue "C40" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D40" '0.131'
Set-TextValue "E40" '  -5.23%  '

Set-TextValue "D41" '0.997'
Set-TextValue "E41" '  -0.32%  '

Set-TextValue "D42" '0.395'
Set-TextValue "E42" '  +0.22%  '

Set-TextValue "E43" '  +0.26%  '

Set-TextValue "E44" '  -0.02%  '

Set-TextValue "D45" '149.90'
Set-TextValue "E45" '  -0.70%  '

Set-TextValue "B46" 'OKB'
Set-TextValue "C46" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D46" '45.57'
Set-TextValue "E46" '  +2.02%  '

Set-TextValue "B47" 'Aave'
Set-TextValue "C47" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D47" '179.91'
Set-TextValue "E47" '  +1.07%  '

Set-TextValue "D48" '0.128'
Set-TextValue "E48" '  +2.11%  '

Set-TextValue "D49" '1.29'
Set-TextValue "E49" '  -0.04%  '

Set-TextValue "B50" 'ARBITRUM'
Set-TextValue "C50" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D50" '0.639'
Set-TextValue "E50" '  +2.89%  '

Set-TextValue "B51" 'Filecoin'
Set-TextValue "C51" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D51" '4.23'
Set-TextValue "E51" '  +0.37%  '
